$wb = $excel.ActiveWorkbook

# Append a new blank worksheet at the end of the workbook, named "08-10-24"
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "08-10-24"

# Header row
$ws.Range("A1").Value = "Games"
$ws.Range("B1").Value = "Score"
$headerRange = $ws.Range("A1:B1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Data rows
$ws.Cells.Item(2, 1).Value = "('BOS', 'HOU')"
$ws.Cells.Item(2, 2).Value = 0.738
$ws.Cells.Item(3, 1).Value = "('ATL', 'COL')"
$ws.Cells.Item(3, 2).Value = 0.731
$ws.Cells.Item(4, 1).Value = "('NYM', 'SEA')"
$ws.Cells.Item(4, 2).Value = 0.721
$ws.Cells.Item(5, 1).Value = "('BAL', 'TB')"
$ws.Cells.Item(5, 2).Value = 0.715
$ws.Cells.Item(6, 1).Value = "('DET', 'SF')"
$ws.Cells.Item(6, 2).Value = 0.698
$ws.Cells.Item(7, 1).Value = "('KC', 'STL')"
$ws.Cells.Item(7, 2).Value = 0.674
$ws.Cells.Item(8, 1).Value = "('CHC', 'CWS')"
$ws.Cells.Item(8, 2).Value = 0.522
$ws.Cells.Item(9, 1).Value = "('OAK', 'TOR')"
$ws.Cells.Item(9, 2).Value = 0.522
$ws.Cells.Item(10, 1).Value = "('LAD', 'PIT')"
$ws.Cells.Item(10, 2).Value = 0.501
$ws.Cells.Item(11, 1).Value = "('AZ', 'PHI')"
$ws.Cells.Item(11, 2).Value = 0.478
$ws.Cells.Item(12, 1).Value = "('CLE', 'MIN')"
$ws.Cells.Item(12, 2).Value = 0.47
$ws.Cells.Item(13, 1).Value = "('LAA', 'WSH')"
$ws.Cells.Item(13, 2).Value = 0.272
$ws.Cells.Item(14, 1).Value = "('NYY', 'TEX')"
$ws.Cells.Item(14, 2).Value = 0.18
$ws.Cells.Item(15, 1).Value = "('CIN', 'MIL')"
$ws.Cells.Item(15, 2).Value = 0.06
$ws.Cells.Item(16, 1).Value = "('MIA', 'SD')"
$ws.Cells.Item(16, 2).Value = -0.025
